$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts D:K -> E:L), to make room for new fiscal-year data.
$ws.Range("D1").EntireColumn.Insert()

# Copy number formatting from column E (the old column D data, now shifted) into the new column D
# so the new column keeps the same look (date format for the header row, number format for data rows).
$ws.Range("E7:E102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new column D with the latest fiscal year figures.
    $ws.Range("D7").Value = "12/31/2018"
    $ws.Range("D8").Value = 169000
    $ws.Range("D9").Value = "NA"
    $ws.Range("D10").Value = "NA"
    $ws.Range("D12").Value = "NA"
    $ws.Range("D13").Value = 0
    $ws.Range("D14").Value = 0
    $ws.Range("D15").Value = -900
    $ws.Range("D17").Value = 34000
    $ws.Range("D18").Value = 135000
    $ws.Range("D20").Value = -86600
    $ws.Range("D21").Value = 53100
    $ws.Range("D23").Value = 48400
    $ws.Range("D24").Value = 9100
    $ws.Range("D26").Value = 39200
    $ws.Range("D27").Value = 38400
    $ws.Range("D28").Value = 0
    $ws.Range("D29").Value = 0
    $ws.Range("D30").Value = 0
    $ws.Range("D31").Value = 0
    $ws.Range("D32").Value = 86600
    $ws.Range("D33").Value = 38400
    $ws.Range("D34").Value = 0
    $ws.Range("D35").Value = 38400
    $ws.Range("D38").Value = "12/31/2018"
    $ws.Range("D41").Value = 142100
    $ws.Range("D42").Value = 177300
    $ws.Range("D43").Value = 0
    $ws.Range("D44").Value = 0
    $ws.Range("D45").Value = 0
    $ws.Range("D46").Value = 0
    $ws.Range("D47").Value = 0
    $ws.Range("D48").Value = 35000
    $ws.Range("D49").Value = 110300
    $ws.Range("D50").Value = 0
    $ws.Range("D51").Value = 0
    $ws.Range("D52").Value = 10300
    $ws.Range("D53").Value = 0
    $ws.Range("D54").Value = 4700700
    $ws.Range("D57").Value = "NA"
    $ws.Range("D58").Value = 0
    $ws.Range("D59").Value = 0
    $ws.Range("D60").Value = 0
    $ws.Range("D61").Value = 78800
    $ws.Range("D62").Value = 0
    $ws.Range("D63").Value = 0
    $ws.Range("D64").Value = 0
    $ws.Range("D65").Value = 0
    $ws.Range("D66").Value = 4246900
    $ws.Range("D68").Value = 0
    $ws.Range("D69").Value = 0
    $ws.Range("D70").Value = 0
    $ws.Range("D71").Value = 0
    $ws.Range("D72").Value = 117400
    $ws.Range("D73").Value = 0
    $ws.Range("D74").Value = 0
    $ws.Range("D75").Value = 0
    $ws.Range("D76").Value = 453800
    $ws.Range("D77").Value = 0
    $ws.Range("D80").Value = "12/31/2018"
    $ws.Range("D81").Value = 38400
    $ws.Range("D83").Value = 4700
    $ws.Range("D84").Value = 0
    $ws.Range("D85").Value = 0
    $ws.Range("D86").Value = 0
    $ws.Range("D87").Value = 0
    $ws.Range("D88").Value = 0
    $ws.Range("D89").Value = 58400
    $ws.Range("D91").Value = -5300
    $ws.Range("D92").Value = 0
    $ws.Range("D93").Value = 0
    $ws.Range("D94").Value = -80500
    $ws.Range("D96").Value = -18300
    $ws.Range("D97").Value = 0
    $ws.Range("D98").Value = 0
    $ws.Range("D99").Value = 0
    $ws.Range("D100").Value = 222800
    $ws.Range("D101").Value = 0
    $ws.Range("D102").Value = 200600
